# Adjust negative test case: on the "pp.xpt" sheet, cell M5 (PPSTRESU,
# "Standard Units") previously held the now-unused value "pmol/L/ug".
# Update it to match the Original Units value "day*ug/mL/mg" (same value
# already present in J5), so the "pmol/L/ug" shared string becomes
# unreferenced and is dropped from the workbook on save.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("pp.xpt")
$ws.Range("M5").Value = "day*ug/mL/mg"
